# Add waveform generator "update rate" register to both analog-output
# blocks (AO0 / AO1) on Sheet1.
#
# Also: rename the first Analog Output block's label from
# "Analog Output" to "AO0", and add a whole second block "AO1" (a
# duplicate of the AO0 register map, with its own UID) right below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the existing "Analog Output" subsystem label to "AO0" ---
$ws.Range("A7").Value = "AO0"

# --- make room for the new "update rate" register inside the AO0 block ---
# Old layout (rows 7-11): idle offset(0), active?(1), increment(2),
# number of elements(3), LUT(4..1026).
# New layout (rows 7-12): idle offset(0), active?(1), update rate(2),
# increment(3), number of elements(4), LUT(5..1026).
$ws.Rows("9:9").Insert()

$ws.Cells.Item(9, 8).Value = 2
$ws.Cells.Item(9, 9).Value = "U32"
$ws.Cells.Item(9, 10).Value = "W"
$ws.Cells.Item(9, 11).Value = "update rate"

$ws.Cells.Item(10, 8).Value = 3
$ws.Cells.Item(10, 9).Value = "U32"
$ws.Cells.Item(10, 10).Value = "W"
$ws.Cells.Item(10, 11).Value = "increment"

$ws.Cells.Item(11, 8).Value = 4
$ws.Cells.Item(11, 9).Value = "U16"
$ws.Cells.Item(11, 10).Value = "W"
$ws.Cells.Item(11, 11).Value = "number of elements"

$ws.Cells.Item(12, 8).Value = "5..1026"
$ws.Cells.Item(12, 9).Value = "I16"
$ws.Cells.Item(12, 10).Value = "W"
$ws.Cells.Item(12, 11).Value = "LUT"

# --- add the new "AO1" subsystem block (rows 14-19), a twin of AO0 ---
$ws.Range("A14").Value = "AO1"
$ws.Cells.Item(14, 2).Value = 32770
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = "U16"
$ws.Cells.Item(14, 10).Value = "W"
$ws.Cells.Item(14, 11).Value = "idle offset"

$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = "BOOL"
$ws.Cells.Item(15, 10).Value = "W"
$ws.Cells.Item(15, 11).Value = "active?"

$ws.Cells.Item(16, 8).Value = 2
$ws.Cells.Item(16, 9).Value = "U32"
$ws.Cells.Item(16, 10).Value = "W"
$ws.Cells.Item(16, 11).Value = "update rate"

$ws.Cells.Item(17, 8).Value = 3
$ws.Cells.Item(17, 9).Value = "U32"
$ws.Cells.Item(17, 10).Value = "W"
$ws.Cells.Item(17, 11).Value = "increment"

$ws.Cells.Item(18, 8).Value = 4
$ws.Cells.Item(18, 9).Value = "U16"
$ws.Cells.Item(18, 10).Value = "W"
$ws.Cells.Item(18, 11).Value = "number of elements"

$ws.Cells.Item(19, 8).Value = "5..1026"
$ws.Cells.Item(19, 9).Value = "I16"
$ws.Cells.Item(19, 10).Value = "W"
$ws.Cells.Item(19, 11).Value = "LUT"

# --- cosmetic: window scroll position + last active cell/selection ---
$wb.Windows.Item(1).Left = 2790
$ws.Range("D22").Select()
